$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift the A:K block up by one row (drops the old sub-header row that
#    lived in row 2 and moves every data/filler row up by one). Columns
#    L:M are intentionally left untouched by using a column-restricted
#    range instead of an EntireRow delete.
$ws.Range("A2:K2").Delete(-4162)

# 2) The shift above drops the very last filler row (91) because its
#    A:K cells are now empty across the whole row. Re-create the
#    original L91:M91 placeholders (still styled, still empty) by
#    copying the format from the row above.
$ws.Range("L90:M90").Copy()
$ws.Range("L91:M91").PasteSpecial(-4122)

# 3) Build the new header style (General number format + the small
#    9pt Arial font used throughout the sheet, i.e. same look as the
#    existing text style but without forcing applyNumberFormat).
$tmpStyle = $wb.Styles.Add("__TmpHeaderStyle")
$tmpStyle.Font.Name = "Arial"
$tmpStyle.Font.Size = 9
$ws.Range("F1:K1").Style = "__TmpHeaderStyle"
$tmpStyle.Delete()

# 4) Write the new single header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# 5) Match the author's final selection.
$ws.Range("A2:K2").Select()
